# Sprint 2 Backlog - "Added Sprint story points"
#
# Story "Display map" (row 5) gets its Day 3 / Day 4 / Day 5 story-point
# estimates filled in. The Day totals in row 34 recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2

# Move the selection/scroll position to A6 (matches the saved view state).
[void]$ws.Range("A6").Select()
